$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only the trailing K2 value is removed; B2:J2 stay as-is
$ws.Range("K2").ClearContents()

# Rows 3-24: every row gets a brand new value inserted at column B,
# with the rest of that row's values shifted one column to the right
# (the former right-most value in each row is dropped).
$rowsData = @(
    ,@(3, @(2.174397541324862, 10.04419022232486, -8.284815777675137, -0.1447337776751375, 1.643575222324863, -1.534543777675137, -1.829872777675137, 0.6845722223248626, -0.2451497776751375, -0.01014177767513746))
    ,@(4, @(8.987815262142332, -9.341190737857668, -1.201108737857667, 0.5872002621423329, -2.590918737857667, -2.886247737857667, -0.371802737857667, -1.301524737857667, -1.066516737857667, -1.373652737857667))
    ,@(5, @(-8.470544361886736, -0.3304623618867365, 1.457846638113264, -1.720272361886736, -2.015601361886736, 0.4988436381132635, -0.4308783618867365, -0.1958703618867365, -0.5030063618867365, 0.007669638113263488))
    ,@(6, @(-0.34111357821662, 1.44719542178338, -1.73092357821662, -2.02625257821662, 0.48819242178338, -0.44152957821662, -0.20652157821662, -0.5136575782166199, -0.002981578216620012, -0.29034857821662))
    ,@(7, @(1.338580889567239, -1.839538110432761, -2.134867110432761, 0.3795778895672392, -0.5501441104327608, -0.3151361104327608, -0.6222721104327609, -0.1115961104327608, -0.3989631104327608, -0.3709641104327608))
    ,@(8, @(-1.631436389465022, -1.926765389465022, 0.5876796105349782, -0.3420423894650218, -0.1070343894650218, -0.4141703894650218, 0.09650561053497819, -0.1908613894650218, -0.1628623894650218, -0.0702143894650218))
    ,@(9, @(-1.810267219091221, 0.7041777809087787, -0.2255442190912214, 0.009463780908778634, -0.2976722190912214, 0.2130037809087786, -0.07436321909122137, -0.04636421909122137, 0.04628378090877863, 0.2592597809087786))
    ,@(10, @(0.7496711949059137, -0.1800508050940863, 0.05495719490591373, -0.2521788050940862, 0.2584971949059137, -0.02886980509408627, -0.0008708050940862688, 0.09177719490591374, 0.3047531949059137, -0.1959068050940863))
    ,@(11, @(-0.1976049264540507, 0.03740307354594928, -0.2697329264540507, 0.2409430735459493, -0.04642392645405072, -0.01842492645405072, 0.07422307354594927, 0.2871990735459493, -0.2134609264540507, 0.07629607354594928))
    ,@(12, @(0.1255103924969555, -0.1816256075030445, 0.3290503924969555, 0.0416833924969555, 0.06968239249695551, 0.1623303924969555, 0.3753063924969555, -0.1253536075030445, 0.1644033924969555, -0.0934886075030445))
    ,@(13, @(-0.003050974015260888, 0.5076250259847391, 0.2202580259847391, 0.2482570259847391, 0.3409050259847391, 0.5538810259847391, 0.0532210259847391, 0.3429780259847391, 0.08508602598473911, 0.4674670259847391))
    ,@(14, @(1.419652253737389, 1.132285253737389, 1.160284253737389, 1.252932253737389, 1.465908253737389, 0.9652482537373891, 1.255005253737389, 0.9971132537373891, 1.379494253737389, 0.9682171597373891))
    ,@(15, @(0.1883110177716323, 0.2163100177716323, 0.3089580177716323, 0.5219340177716323, 0.02127401777163229, 0.3110310177716323, 0.0531390177716323, 0.4355200177716323, 0.02424292377163231, 0.3111990177716323))
    ,@(16, @(0.2202779152847414, 0.3129259152847414, 0.5259019152847414, 0.0252419152847414, 0.3149989152847414, 0.05710691528474141, 0.4394879152847414, 0.02821082128474142, 0.3151669152847414))
    ,@(17, @(0.4485660054549828, 0.6615420054549828, 0.1608820054549828, 0.4506390054549828, 0.1927470054549828, 0.5751280054549828, 0.1638509114549828, 0.4508070054549828))
    ,@(18, @(0.420735823599318, -0.07992417640068192, 0.2098328235993181, -0.04805917640068191, 0.3343218235993181, -0.0769552704006819, 0.2100008235993181))
    ,@(19, @(-0.1252553916527783, 0.1645016083472217, -0.09339039165277828, 0.2889906083472217, -0.1222864856527783, 0.1646696083472217))
    ,@(20, @(0.09764018641116785, -0.1602518135888321, 0.2221291864111679, -0.1891479075888321, 0.09780818641116784))
    ,@(21, @(-0.1898380159455487, 0.1925429840544512, -0.2187341099455487, 0.06822198405445125))
    ,@(22, @(0.1743923273248104, -0.2368847666751896, 0.05007132732481036))
    ,@(23, @(-0.254916590923889, 0.03203950307611103))
    ,@(24, @(-0.01954473626955232))
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
